$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Duplicate the "总计" sheet to create "2022-Q1", positioned
#    right before "总计" (i.e. right after "2021-Q2"). Copying keeps
#    the sheet's pagesetup/format scaffolding (sheetPr, pageMargins,
#    header/"A" column styling) intact.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("总计 (2)")
$newSheet.Name = "2022-Q1"

# Extend the styled header range (B1:D1) across the new E1:H1 columns
# by copying the formatting of an existing styled header cell first.
$newSheet.Range("D1").Copy($newSheet.Range("E1:H1"))

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row (row 2). A2 already carries the right style from the copy.
$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "519929"
$newSheet.Range("C2").Value = "长信电子信息行业量化灵活配置混合"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.00"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "93.92"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "4.34"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0434"
$newSheet.Range("H2").Value = 9

# ------------------------------------------------------------------
# 2. Insert a new row 2 in "总计" sheet for the "2022-Q1" summary,
#    pushing the existing "2021-Q2" row down to row 3
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("总计")
$ws2.Rows.Item(2).Insert()
$ws2.Range("B2:D2").ClearFormats()

# Give the new row's "A" cell the same styling as the sibling "A"
# cells in this column by copying the format straight from A3 (the
# shifted-down original row, which already carries the correct style).
$ws2.Range("A3").Copy($ws2.Range("A2"))
$ws2.Range("A2").Value = 0

$ws2.Range("B2").Value = "2022-Q1"
$ws2.Range("C2").Value = 1
$ws2.Range("D2").Value = 0.04

# The shifted-down original row's running index (column A) advances
# from 0 to 1 now that it is the second data row.
$ws2.Range("A3").Value = 1

# ------------------------------------------------------------------
# 3. Restore the originally active sheet/tab
# ------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
